# Add 2022-Q3 data: insert a new sheet "2022-Q3" right after "总计",
# and add a corresponding summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for "2022-Q3"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# New row 2 content/style
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = "4.78"

# Renumber the (now shifted down) rows 3..9 sequentially
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: clone "2022-Q1" (same 5-fund shape) and
#    overwrite with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Keep fund-code / numeric-looking-text columns as TEXT (leading zeros,
# "27.40" etc. must not turn into numbers).
$q3.Range("B2:B6").NumberFormat = "@"
$q3.Range("D2:G5").NumberFormat = "@"

# Row 2 - 001838
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "001838"
$q3.Range("C2").Value = "国投瑞银国家安全灵活配置混合"
$q3.Range("D2").Value = "27.40"
$q3.Range("E2").Value = "94.42"
$q3.Range("F2").Value = "7.21"
$q3.Range("G2").Value = "1.9755"
$q3.Range("H2").Value = 9

# Row 3 - 005774
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "005774"
$q3.Range("C3").Value = "华夏产业升级混合A"
$q3.Range("D3").Value = "24.29"
$q3.Range("E3").Value = "93.85"
$q3.Range("F3").Value = "7.10"
$q3.Range("G3").Value = "1.7246"
$q3.Range("H3").Value = 5

# Row 4 - 015059
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "015059"
$q3.Range("C4").Value = "华夏产业升级混合C"
$q3.Range("D4").Value = "8.92"
$q3.Range("E4").Value = "93.85"
$q3.Range("F4").Value = "7.10"
$q3.Range("G4").Value = "0.6333"
$q3.Range("H4").Value = 5

# Row 5 - 460002
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "460002"
$q3.Range("C5").Value = "华泰柏瑞积极成长混合A"
$q3.Range("D5").Value = "5.83"
$q3.Range("E5").Value = "90.49"
$q3.Range("F5").Value = "7.71"
$q3.Range("G5").Value = "0.4495"
$q3.Range("H5").Value = 3

# Row 6 - 960030 (market value is genuinely 0, stored as a number)
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "960030"
$q3.Range("C6").Value = "华泰柏瑞积极成长混合H"
$q3.Range("D6").NumberFormat = "@"
$q3.Range("D6").Value = "0.00"
$q3.Range("E6").NumberFormat = "@"
$q3.Range("E6").Value = "90.49"
$q3.Range("F6").NumberFormat = "@"
$q3.Range("F6").Value = "7.71"
$q3.Range("G6").Value = 0
$q3.Range("H6").Value = 3

# ---------------------------------------------------------------------
# 3) Restore the originally active tab ("2020-Q4", now shifted to the
#    last position).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "2022-Q3 sheet added"
